$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# "Experimental" row (row 7): set the Value cell (B7) to the literal text "true"
# A bare Value = "true" gets auto-coerced to a Boolean by Excel, so a leading
# apostrophe is used to force text entry; the apostrophe also flags the cell's
# style with a "quote prefix" marker, so re-apply B6's (identical) formatting
# afterwards to restore the original, unmarked cell style.
$ws.Range("B7").Value = "'true"
$ws.Range("B6").Copy() | Out-Null
$ws.Range("B7").PasteSpecial(-4122) | Out-Null

# "Date" row (row 8): update the Value cell (B8) with the new timestamp
$ws.Range("B8").Value = "2023-02-16T14:43:10-06:00"
